# Add two new columns ("min" and "max") right after the "right answer"
# column (column G) and before the existing "question_type" / "solution"
# columns, pushing the latter two from H:I to J:K.
#
# Resulting header row:
#   A: question       B: option A   C: option B   D: option C
#   E: option D       F: option E   G: right answer
#   H: min (new)      I: max (new)
#   J: question_type  K: solution

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at H:I; this shifts the current H ("question_type")
# and I ("solution") columns to J and K respectively, and extends the sheet's
# used range from A1:I1 to A1:K1.
$ws.Columns("H:I").Insert()

# Fill the two freshly-inserted header cells.
$ws.Range("H1").Value = "min"
$ws.Range("I1").Value = "max"

# Match the workbook's recorded selection/active cell after the edit.
$ws.Range("H2").Select()
